$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-19 Tuesday" "2024-11-20 Wednesday"
Replace-Text "135÷7=19, 2" "944÷4=236, 0"
Replace-Text "452÷3=150, 2" "540÷4=135, 0"
Replace-Text "910÷3=303, 1" "262÷8=32, 6"
Replace-Text "101÷2=50, 1" "957÷6=159, 3"
Replace-Text "823÷3=274, 1" "250÷9=27, 7"
Replace-Text "544÷6=90, 4" "542÷8=67, 6"
Replace-Text "625÷9=69, 4" "278÷3=92, 2"
Replace-Text "458÷5=91, 3" "216÷6=36, 0"
Replace-Text "356÷4=89, 0" "497÷3=165, 2"
Replace-Text "895÷5=179, 0" "692÷3=230, 2"
Replace-Text "741÷3=247, 0" "501÷8=62, 5"
Replace-Text "766÷9=85, 1" "614÷9=68, 2"
Replace-Text "778÷8=97, 2" "449÷7=64, 1"
Replace-Text "143÷9=15, 8" "374÷4=93, 2"
Replace-Text "188÷2=94, 0" "116÷4=29, 0"
Replace-Text "242÷2=121, 0" "981÷3=327, 0"
Replace-Text "724÷9=80, 4" "181÷3=60, 1"
Replace-Text "510÷5=102, 0" "933÷7=133, 2"
Replace-Text "512÷6=85, 2" "242÷4=60, 2"
Replace-Text "496÷4=124, 0" "920÷5=184, 0"
Replace-Text "793÷3=264, 1" "519÷4=129, 3"
Replace-Text "193÷8=24, 1" "859÷5=171, 4"
Replace-Text "842÷7=120, 2" "150÷2=75, 0"
Replace-Text "912÷8=114, 0" "591÷2=295, 1"
Replace-Text "385÷2=192, 1" "601÷6=100, 1"
